$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("G3").Value = 1.91
$ws.Range("J3").Value = 2.88
$ws.Range("L3").Value = 4.75
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("AR3").Value = 1.82
$ws.Range("AS3").Value = 1.92

# Row 7 updates
$ws.Range("Q7").Value = 1.58
$ws.Range("R7").Value = 2.25
$ws.Range("S7").Value = 2.05
$ws.Range("T7").Value = 1.8
